$d = $word.ActiveDocument

# NOTE on ordering: the "gradient descent" block (originally slope 1.0125,
# coefficients 2.9906/1.9762) is updated to share the same slope value
# (0.965) and the same beta1 value (2.9906) that the "scikit-learn" and
# "custom OLS" blocks are updated *to*. To avoid a later global
# find/replace re-matching text that an earlier replacement just produced,
# we perform the gradient-descent block's replacements FIRST (while its
# original values are still unique in the document), and only afterwards
# touch the scikit-learn / custom OLS blocks.

# --- Gradient descent block (bookmark "fit-the-data-using-gradient-descent") ---
$d.Content.Find.Execute("1.0125 and coefficients", $false, $false, $false, $false, $false, $true, 1, $false, "0.965 and coefficients", 2)
$d.Content.Find.Execute("2.9906, and", $false, $false, $false, $false, $false, $true, 1, $false, "2.9894, and", 2)
$d.Content.Find.Execute("1.9762", $false, $false, $false, $false, $false, $true, 1, $false, "2.0108", 2)

# --- scikit-learn block + custom OLS block (both identical, updated together) ---
$d.Content.Find.Execute("1.0132 and coefficients", $false, $false, $false, $false, $false, $true, 1, $false, "0.965 and coefficients", 2)
$d.Content.Find.Execute("2.9918, and", $false, $false, $false, $false, $false, $true, 1, $false, "2.9906, and", 2)
$d.Content.Find.Execute("1.9769", $false, $false, $false, $false, $false, $true, 1, $false, "2.0115", 2)

# --- numpy array() formatted block near the end of the document ---
$d.Content.Find.Execute("array([1.01726883]) and coefficients", $false, $false, $false, $false, $false, $true, 1, $false, "array([0.94343469]) and coefficients", 2)
$d.Content.Find.Execute("array([2.99623763]), and", $false, $false, $false, $false, $false, $true, 1, $false, "array([2.98563073]), and", 2)
$d.Content.Find.Execute("array([1.94734374])", $false, $false, $false, $false, $false, $true, 1, $false, "array([2.04592354])", 2)
